$d = $word.ActiveDocument

# Helper: find the paragraph whose text contains $needle.
function Find-ParagraphContaining([string]$needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

$semaineParagraph    = Find-ParagraphContaining "3me validation"
$chaqueParagraph     = Find-ParagraphContaining "Chaque validation est propre"
$illustrationsParagraph = Find-ParagraphContaining "illustrations pour mieux montrer"

# ---------------------------------------------------------------------------
# 1) Move the hidden "_GoBack" bookmark from the end of the paragraph ending
#    in "...3me validation : 1 semaine" to the end of the paragraph
#    "- Chaque validation est propre a l'utilisateur (...)".
#
#    A collapsed (zero-length) Range placed exactly at a paragraph's last
#    character position trips an edge case in Bookmarks.Add, so we briefly
#    insert a one-character placeholder at that spot, bookmark the
#    placeholder range, then delete the placeholder again - the bookmark
#    survives, collapsed at the original location, which is what we want.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$target = $chaqueParagraph.Range.Duplicate
$target.Collapse(0)          # collapse to just after the paragraph's last character
[void]$target.MoveStart(1, -1)
[void]$target.MoveEnd(1, -1) # $target is now the collapsed point right before the pilcrow
$markerStart = $target.Start
$target.InsertAfter([char]0xE000)   # temporary 1-char placeholder

$markerRange = $d.Range($markerStart, $markerStart + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$d.Range($markerStart, $markerStart + 1).Delete()   # remove placeholder; bookmark stays put

# ---------------------------------------------------------------------------
# 2) Remove the paragraph about illustrations entirely.
# ---------------------------------------------------------------------------
$illustrationsParagraph.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Drop one of the trailing empty paragraphs (4 -> 3).
# ---------------------------------------------------------------------------
$firstEmptyAfter = $null
$passedChaque = $false
foreach ($p in $d.Paragraphs) {
    if ($passedChaque -and $p.Range.Text.Trim().Length -eq 0) {
        $firstEmptyAfter = $p
        break
    }
    if ($p.Range.Text -like "*Chaque validation est propre*") {
        $passedChaque = $true
    }
}
$firstEmptyAfter.Range.Delete()
